$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each pair of rows below had their full record (columns B..AD) swapped between
# the two rows, while the leading row-index column (A) stayed fixed in place.
$rowPairs = @(
    @(18, 19),
    @(43, 45),
    @(52, 53),
    @(82, 83),
    @(88, 89),
    @(114, 115),
    @(120, 121),
    @(124, 126),
    @(146, 147),
    @(174, 175),
    @(206, 207),
    @(208, 209),
    @(216, 217),
    @(218, 219),
    @(228, 229),
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $rng1 = $ws.Range($ws.Cells.Item($r1, 2), $ws.Cells.Item($r1, 30))
    $rng2 = $ws.Range($ws.Cells.Item($r2, 2), $ws.Cells.Item($r2, 30))
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value2 = $v2
    $rng2.Value2 = $v1
}

Write-Output "swap complete"